$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet held a 3-column header (firstName/lastName/handle) in row 1.
# The new sheet is a single "numbers" column, so drop the now-unused B/C header cells.
$ws.Cells.Item(1, 2).ClearContents()
$ws.Cells.Item(1, 3).ClearContents()

# New header
$ws.Range("A1").Value = "numbers"

# Seed values: A2 is a literal, A3 is a lone formula, A4:A30 is one fill-down block.
$ws.Range("A2").Value = 1098
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A30").Formula = "=A3+1"

# Row 31 was retyped by hand in the source workbook (breaking the fill block in two),
# so enter it as its own standalone formula rather than part of a filled range.
$ws.Range("A31").Formula = "=A30+1"

# A32:A43 is the second fill-down block, continuing on from the manually entered A31.
$ws.Range("A32:A43").Formula = "=A31+1"

$ws.Range("B31").Select()
